# Applies scheduled market-data refresh to the Leve profit sheets.
# For each touched row, columns H:N (currentAveragePrice.., LevePriceNQ/HQ,
# LeveProfitNQ/HQ) are refreshed to the latest computed values. A few rows
# gain or lose an H:N cell entirely where a profit figure newly exists or
# no longer applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @(12, 8, 267.85715),
    @(12, 9, 250),
    @(12, 10, 312.5),
    @(12, 11, 250),
    @(12, 12, 312.5),
    @(12, 13, -80),
    @(12, 14, -652.5),
    @(20, 8, 883.5),
    @(20, 9, 883.5),
    @(20, 11, 883.5),
    @(20, 13, -653.5),
    @(35, 8, 883.5),
    @(35, 9, 883.5),
    @(35, 11, 883.5),
    @(35, 13, -504.5),
    @(38, 8, 881.9231),
    @(38, 9, 58.125),
    @(38, 10, 2200),
    @(38, 11, 174.375),
    @(38, 12, 6600),
    @(38, 13, 197.625),
    @(38, 14, -7344),
    @(58, 9, 687.5),
    @(58, 10, 1901.4),
    @(58, 11, 2062.5),
    @(58, 12, 5704.200000000001),
    @(58, 13, -1912.5),
    @(58, 14, -6004.200000000001),
    @(70, 8, 4498.727),
    @(70, 10, 4498.727),
    @(70, 12, 13496.181),
    @(70, 14, -14036.181),
    @(73, 8, 4498.727),
    @(73, 10, 4498.727),
    @(73, 12, 13496.181),
    @(73, 14, -15368.181),
    @(80, 8, 854.1667),
    @(80, 10, 474.75),
    @(80, 12, 1424.25),
    @(80, 14, -3420.25),
    @(83, 8, 854.1667),
    @(83, 10, 474.75),
    @(83, 12, 4272.75),
    @(83, 14, -14256.75),
    @(87, 8, 40000),
    @(87, 10, 40000),
    @(87, 12, 40000),
    @(87, 14, -42496),
    @(88, 8, 894.9167),
    @(88, 9, 248.8),
    @(88, 10, 1356.4286),
    @(88, 11, 248.8),
    @(88, 12, 1356.4286),
    @(88, 13, 157.2),
    @(88, 14, -2168.4286),
    @(90, 8, 40000),
    @(90, 10, 40000),
    @(90, 12, 120000),
    @(90, 14, -132480),
    @(91, 8, 894.9167),
    @(91, 9, 248.8),
    @(91, 10, 1356.4286),
    @(91, 11, 248.8),
    @(91, 12, 1356.4286),
    @(91, 13, 1155.2),
    @(91, 14, -4164.4286),
    @(125, 8, 7599.6),
    @(125, 9, 4666.3335),
    @(125, 10, 11999.5),
    @(125, 11, 41997.0015),
    @(125, 12, 107995.5),
    @(125, 13, -39537.0015),
    @(125, 14, -112915.5),
    @(138, 8, 1913.4546),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @(74, 8, 1490.8572),
    @(74, 9, 1498.3334),
    @(74, 11, 1498.3334),
    @(74, 13, -624.3334),
    @(77, 8, 1490.8572),
    @(77, 9, 1498.3334),
    @(77, 11, 7491.666999999999),
    @(77, 13, -3123.666999999999),
    @(88, 8, 1844.4166),
    @(88, 10, 2371.4119),
    @(88, 12, 2371.4119),
    @(88, 14, -3183.4119),
    @(91, 8, 1844.4166),
    @(91, 10, 2371.4119),
    @(91, 12, 2371.4119),
    @(91, 14, -5179.4119),
    @(92, 8, 59500),
    @(92, 10, 59500),
    @(92, 12, 59500),
    @(92, 14, -64492),
    @(122, 8, 1895),
    @(122, 9, 1841.1111),
    @(122, 11, 5523.3333),
    @(122, 13, -3073.3333),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @(86, 8, 1883.5834),
    @(86, 9, 1134.3684),
    @(86, 10, 4730.6),
    @(86, 11, 1134.3684),
    @(86, 12, 4730.6),
    @(86, 13, -11.36840000000007),
    @(86, 14, -6976.6),
    @(89, 8, 1883.5834),
    @(89, 9, 1134.3684),
    @(89, 10, 4730.6),
    @(89, 11, 5671.842000000001),
    @(89, 12, 23653),
    @(89, 13, -55.84200000000055),
    @(89, 14, -34885),
    @(94, 8, 942.6799999999999),
    @(94, 9, 915.2222),
    @(94, 10, 1013.2857),
    @(94, 11, 915.2222),
    @(94, 12, 1013.2857),
    @(94, 13, -464.2222),
    @(94, 14, -1915.2857),
    @(99, 8, 3077.5557),
    @(99, 9, 3099.7144),
    @(99, 11, 3099.7144),
    @(99, 13, -1601.7144),
    @(105, 8, 1969.8),
    @(105, 9, 2099.5),
    @(105, 11, 2099.5),
    @(105, 13, -352.5),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @(22, 8, 179.33333),
    @(22, 9, 175.2),
    @(22, 11, 175.2),
    @(22, 13, 174.8),
    @(105, 8, 703.6),
    @(105, 9, 634.4375),
    @(105, 11, 634.4375),
    @(105, 13, 1112.5625),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @(20, 8, 12150),
    @(20, 10, 29000),
    @(20, 12, 29000),
    @(20, 14, -29490),
    @(24, 8, 0),
    @(24, 10, 0),
    @(24, 12, 0),
    @(24, 14, $null),
    @(70, 8, 2466.3333),
    @(70, 9, 2499),
    @(70, 10, 2450),
    @(70, 11, 2499),
    @(70, 12, 2450),
    @(70, 13, -2229),
    @(70, 14, -2990),
    @(73, 8, 2466.3333),
    @(73, 9, 2499),
    @(73, 10, 2450),
    @(73, 11, 2499),
    @(73, 12, 2450),
    @(73, 13, -1563),
    @(73, 14, -4322),
    @(126, 8, 2784.1428),
    @(126, 9, 2748.1667),
    @(126, 11, 8244.500100000001),
    @(126, 13, -5774.500100000001),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @(7, 8, 4666.6665),
    @(7, 9, 3000),
    @(7, 11, 3000),
    @(7, 13, -2888),
    @(16, 8, 0),
    @(16, 9, 0),
    @(16, 10, 0),
    @(16, 11, 0),
    @(16, 12, 0),
    @(16, 13, $null),
    @(16, 14, $null),
    @(46, 8, 2642.4285),
    @(46, 9, 2642.4285),
    @(46, 11, 2642.4285),
    @(46, 13, -2454.4285),
    @(93, 8, 432.66666),
    @(93, 9, 432.66666),
    @(93, 11, 432.66666),
    @(93, 13, 815.33334),
    @(100, 8, 2000),
    @(100, 9, 1000),
    @(100, 10, 3000),
    @(100, 11, 1000),
    @(100, 12, 3000),
    @(100, 13, -459),
    @(100, 14, -4082),
    @(126, 8, 4666.6665),
    @(126, 9, 3000),
    @(126, 11, 9000),
    @(126, 13, -6530),
    @(136, 8, 3504),
    @(136, 9, 3504),
    @(136, 11, 10512),
    @(136, 13, -7962),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @(15, 8, 3335335.2),
    @(15, 10, 3003.5),
    @(15, 12, 3003.5),
    @(15, 14, -3579.5),
    @(100, 8, 606),
    @(100, 9, 430.2),
    @(100, 10, 899),
    @(100, 11, 860.4),
    @(100, 12, 1798),
    @(100, 13, -319.4),
    @(100, 14, -2880),
    @(126, 8, 2727.875),
    @(126, 9, 2287.1667),
    @(126, 11, 6861.500100000001),
    @(126, 13, -4391.500100000001),
)
foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e[0], $e[1])
    if ($null -eq $e[2]) {
        $cell.ClearContents()
    } else {
        $cell.Value = $e[2]
    }
}

